# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-123) from 2023-09-13 (45182) to 2023-09-15 (45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C123").Value = 45184
